# Update LR-pair TPM-derived values on the active worksheet.
# This mirrors a re-run of the NATMI scoring pipeline with new TPM input:
# the "Ligand-expressing cells" counts (and therefore every value derived
# from them) changed for several rows, so we overwrite each affected cell
# with its recomputed value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06440233333333333
$ws.Range("H2").Value = 0.193207
$ws.Range("I2").Value = 0.03647206354366116
$ws.Range("J2").Value = 0.03647206354366116
$ws.Range("O2").Value = 0.8416031693647025
$ws.Range("P2").Value = 0.8416031693647025
$ws.Range("Q2").Value = 0.10171575722
$ws.Range("R2").Value = 0.9154418149799999
$ws.Range("S2").Value = 0.03069500427161605
$ws.Range("T2").Value = 0.03069500427161605

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06440233333333333
$ws.Range("H3").Value = 0.193207
$ws.Range("I3").Value = 0.03647206354366116
$ws.Range("J3").Value = 0.03647206354366116
$ws.Range("M3").Value = 0.2972526666666667
$ws.Range("N3").Value = 0.8917580000000001
$ws.Range("O3").Value = 0.1583968306352975
$ws.Range("P3").Value = 0.1583968306352975
$ws.Range("Q3").Value = 0.01914376532288889
$ws.Range("R3").Value = 0.172293887906
$ws.Range("S3").Value = 0.005777059272045105
$ws.Range("T3").Value = 0.005777059272045105

# Row 4
$ws.Range("I4").Value = 0.8194013021867156
$ws.Range("J4").Value = 0.8194013021867155
$ws.Range("O4").Value = 0.8416031693647025
$ws.Range("P4").Value = 0.8416031693647025
$ws.Range("S4").Value = 0.6896107329019041
$ws.Range("T4").Value = 0.6896107329019041

# Row 5
$ws.Range("I5").Value = 0.8194013021867156
$ws.Range("J5").Value = 0.8194013021867155
$ws.Range("M5").Value = 0.2972526666666667
$ws.Range("N5").Value = 0.8917580000000001
$ws.Range("O5").Value = 0.1583968306352975
$ws.Range("P5").Value = 0.1583968306352975
$ws.Range("Q5").Value = 0.4300942888946667
$ws.Range("R5").Value = 3.870848600052
$ws.Range("S5").Value = 0.1297905692848114
$ws.Range("T5").Value = 0.1297905692848114

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2544986666666667
$ws.Range("H6").Value = 0.763496
$ws.Range("I6").Value = 0.1441266342696234
$ws.Range("J6").Value = 0.1441266342696234
$ws.Range("O6").Value = 0.8416031693647025
$ws.Range("P6").Value = 0.8416031693647025
$ws.Range("Q6").Value = 0.4019501041599999
$ws.Range("R6").Value = 3.617550937439999
$ws.Range("S6").Value = 0.1212974321911824
$ws.Range("T6").Value = 0.1212974321911824

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.2544986666666667
$ws.Range("H7").Value = 0.763496
$ws.Range("I7").Value = 0.1441266342696234
$ws.Range("J7").Value = 0.1441266342696234
$ws.Range("M7").Value = 0.2972526666666667
$ws.Range("N7").Value = 0.8917580000000001
$ws.Range("O7").Value = 0.1583968306352975
$ws.Range("P7").Value = 0.1583968306352975
$ws.Range("Q7").Value = 0.07565040732977778
$ws.Range("R7").Value = 0.680853665968
$ws.Range("S7").Value = 0.02282920207844099
$ws.Range("T7").Value = 0.02282920207844099
